# Apply the odds updates described by the XML diff.
# Each block corresponds to one match (row) in the sheet; only the
# specific odds columns that changed between the two workbook snapshots
# are touched here.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 76 Igdir Belediyespor vs Sakaryaspor
$ws.Cells.Item(2, 9).Value = 13.5   # I2: 15 -> 13.5
$ws.Cells.Item(2, 10).Value = 4.4   # J2: 4.3 -> 4.4
$ws.Cells.Item(2, 12).Value = 1.01   # L2: 0 -> 1.01
$ws.Cells.Item(2, 13).Value = 1.01   # M2: 0 -> 1.01
$ws.Cells.Item(2, 14).Value = 1.01   # N2: 0 -> 1.01
$ws.Cells.Item(2, 15).Value = 1.22   # O2: 0 -> 1.22
$ws.Cells.Item(2, 17).Value = 1.22   # Q2: 1.54 -> 1.22
$ws.Cells.Item(2, 18).Value = 1.41   # R2: 0 -> 1.41
$ws.Cells.Item(2, 19).Value = 2.28   # S2: 0 -> 2.28
$ws.Cells.Item(2, 20).Value = 1.01   # T2: 0 -> 1.01
$ws.Cells.Item(2, 21).Value = 1.01   # U2: 0 -> 1.01
$ws.Cells.Item(2, 22).Value = 1.08   # V2: 0 -> 1.08
$ws.Cells.Item(2, 23).Value = 2.42   # W2: 0 -> 2.42
$ws.Cells.Item(2, 24).Value = 1000   # X2: 0 -> 1000
$ws.Cells.Item(2, 25).Value = 1000   # Y2: 0 -> 1000
$ws.Cells.Item(2, 26).Value = 1000   # Z2: 0 -> 1000
$ws.Cells.Item(2, 27).Value = 1000   # AA2: 0 -> 1000
$ws.Cells.Item(2, 28).Value = 1000   # AB2: 0 -> 1000
$ws.Cells.Item(2, 29).Value = 1000   # AC2: 0 -> 1000
$ws.Cells.Item(2, 30).Value = 1000   # AD2: 0 -> 1000
$ws.Cells.Item(2, 31).Value = 1000   # AE2: 0 -> 1000
$ws.Cells.Item(2, 32).Value = 1000   # AF2: 0 -> 1000
$ws.Cells.Item(2, 33).Value = 1000   # AG2: 0 -> 1000
$ws.Cells.Item(2, 34).Value = 1000   # AH2: 0 -> 1000
$ws.Cells.Item(2, 35).Value = 1000   # AI2: 0 -> 1000
$ws.Cells.Item(2, 36).Value = 1000   # AJ2: 0 -> 1000
$ws.Cells.Item(2, 37).Value = 1000   # AK2: 0 -> 1000
$ws.Cells.Item(2, 38).Value = 1000   # AL2: 0 -> 1000
$ws.Cells.Item(2, 39).Value = 1000   # AM2: 0 -> 1000
$ws.Cells.Item(2, 40).Value = 1000   # AN2: 0 -> 1000
$ws.Cells.Item(2, 41).Value = 1000   # AO2: 0 -> 1000

# Row 3: Istanbulspor vs Keciorengucu
$ws.Cells.Item(3, 16).Value = 2.08   # P3: 2.02 -> 2.08
$ws.Cells.Item(3, 17).Value = 1.58   # Q3: 1.57 -> 1.58

# Row 4: Konyaspor vs Eyupspor
$ws.Cells.Item(4, 6).Value = 1.59   # F4: 1.58 -> 1.59

# Row 5: Panaitolikos vs Levadiakos
$ws.Cells.Item(5, 6).Value = 4.7   # F5: 4.9 -> 4.7
$ws.Cells.Item(5, 9).Value = 1.93   # I5: 1.92 -> 1.93
$ws.Cells.Item(5, 22).Value = 2.06   # V5: 2.08 -> 2.06

# Row 6: Omonia FC Aradippou vs Enosis Neon Paralimni
$ws.Cells.Item(6, 17).Value = 2.1   # Q6: 2.08 -> 2.1

# Row 8: Besiktas vs Kayserispor
$ws.Cells.Item(8, 7).Value = 1.46   # G8: 1.47 -> 1.46
$ws.Cells.Item(8, 9).Value = 8.4   # I8: 8.6 -> 8.4
$ws.Cells.Item(8, 16).Value = 2.86   # P8: 2.88 -> 2.86

# Row 10: US Cremonese vs Verona
$ws.Cells.Item(10, 6).Value = 2.68   # F10: 2.66 -> 2.68
$ws.Cells.Item(10, 16).Value = 1.62   # P10: 1.63 -> 1.62

# Row 14: Hapoel Tel Aviv vs Hapoel Beer Sheva
$ws.Cells.Item(14, 16).Value = 1.98   # P14: 1.99 -> 1.98
$ws.Cells.Item(14, 17).Value = 1.01   # Q14: 1.6 -> 1.01

# Row 15: NFC Volos vs Atromitos
$ws.Cells.Item(15, 7).Value = 2.52   # G15: 2.6 -> 2.52
$ws.Cells.Item(15, 9).Value = 5.1   # I15: 5.4 -> 5.1
$ws.Cells.Item(15, 10).Value = 2.94   # J15: 2.84 -> 2.94
$ws.Cells.Item(15, 17).Value = 1.01   # Q15: 2.36 -> 1.01

# Row 16: Jong FC Utrecht vs Jong PSV Eindhoven
$ws.Cells.Item(16, 16).Value = 1.25   # P16: 2.76 -> 1.25
$ws.Cells.Item(16, 17).Value = 1.38   # Q16: 1.31 -> 1.38

# Row 17: Jong Ajax Amsterdam vs Roda JC
$ws.Cells.Item(17, 10).Value = 3.7   # J17: 3.95 -> 3.7
$ws.Cells.Item(17, 16).Value = 2.52   # P17: 2.5 -> 2.52

# Row 19: Granada vs Eibar
$ws.Cells.Item(19, 16).Value = 1.56   # P19: 1.57 -> 1.56

# Row 20: Nancy vs Guingamp
$ws.Cells.Item(20, 6).Value = 2.46   # F20: 2.6 -> 2.46
$ws.Cells.Item(20, 7).Value = 3.3   # G20: 2.96 -> 3.3
$ws.Cells.Item(20, 9).Value = 3   # I20: 3.2 -> 3
$ws.Cells.Item(20, 10).Value = 2.92   # J20: 3.15 -> 2.92
$ws.Cells.Item(20, 11).Value = 3.8   # K20: 3.65 -> 3.8

# Row 21: Lazio vs Como
$ws.Cells.Item(21, 6).Value = 2.92   # F21: 2.9 -> 2.92

# Row 22: Brighton vs Bournemouth
$ws.Cells.Item(22, 7).Value = 1.95   # G22: 1.96 -> 1.95
$ws.Cells.Item(22, 9).Value = 4.2   # I22: 4.1 -> 4.2
$ws.Cells.Item(22, 18).Value = 1.72   # R22: 1.71 -> 1.72
$ws.Cells.Item(22, 20).Value = 1.53   # T22: 1.54 -> 1.53
$ws.Cells.Item(22, 36).Value = 22   # AJ22: 23 -> 22

# Row 23: Elche vs Sevilla
$ws.Cells.Item(23, 6).Value = 2.3   # F23: 2.36 -> 2.3
$ws.Cells.Item(23, 7).Value = 2.34   # G23: 2.4 -> 2.34
$ws.Cells.Item(23, 8).Value = 3.5   # H23: 3.4 -> 3.5
$ws.Cells.Item(23, 9).Value = 3.6   # I23: 3.5 -> 3.6
$ws.Cells.Item(23, 26).Value = 24   # Z23: 27 -> 24
$ws.Cells.Item(23, 32).Value = 14   # AF23: 14.5 -> 14
$ws.Cells.Item(23, 36).Value = 32   # AJ23: 38 -> 32
$ws.Cells.Item(23, 37).Value = 28   # AK23: 29 -> 28

# Row 24: Club Football Estrela vs Estoril Praia
$ws.Cells.Item(24, 16).Value = 1.92   # P24: 1.93 -> 1.92
$ws.Cells.Item(24, 20).Value = 1.72   # T24: 1.59 -> 1.72
$ws.Cells.Item(24, 21).Value = 2.3   # U24: 2.2 -> 2.3
$ws.Cells.Item(24, 29).Value = 9.2   # AC24: 9 -> 9.2

# Row 25: Cucuta Deportivo vs Once Caldas
$ws.Cells.Item(25, 7).Value = 2.5   # G25: 2.54 -> 2.5
$ws.Cells.Item(25, 8).Value = 3.5   # H25: 3.55 -> 3.5
$ws.Cells.Item(25, 9).Value = 4.2   # I25: 4.3 -> 4.2
$ws.Cells.Item(25, 10).Value = 3.05   # J25: 3.1 -> 3.05
$ws.Cells.Item(25, 11).Value = 3.55   # K25: 3.6 -> 3.55
$ws.Cells.Item(25, 17).Value = 2.26   # Q25: 2.22 -> 2.26
$ws.Cells.Item(25, 19).Value = 3.85   # S25: 4.2 -> 3.85
$ws.Cells.Item(25, 23).Value = 1.66   # W25: 1.69 -> 1.66
$ws.Cells.Item(25, 25).Value = 14   # Y25: 14.5 -> 14
